$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.391.40'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '2.920.27'
$ws.Range("E3").Value = '  -2.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '375.67'
$ws.Range("E5").Value = '  +6.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.57'
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("E7").Value = '  -2.79%  '
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  -4.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.92'
$ws.Range("E10").Value = '  -2.93%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.28'
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").Value = '3.381.03'
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").Value = '2.913.30'
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("E17").Value = '  -8.66%  '
$ws.Range("D18").Value = '51.302.26'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.38'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.31'
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.90'
$ws.Range("E21").Value = '  -4.51%  '
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.27'
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.43'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("E26").Value = '  -4.51%  '
$ws.Range("E27").Value = '  -5.02%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.62'
$ws.Range("E29").Value = '  -4.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.84'
$ws.Range("E31").Value = '  +7.87%  '
$ws.Range("E32").Value = '  -4.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.79'
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.12'
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  -5.53%  '
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("E39").Value = '  -10.51%  '
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("E41").Value = '  -10.52%  '
$ws.Range("E42").Value = '  -7.64%  '
$ws.Range("E43").Value = '  -2.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.51'
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.56'
$ws.Range("E45").Value = '  -5.22%  '
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.270'
$ws.Range("E47").Value = '  +11.26%  '
$ws.Range("D48").Value = '2.021.89'
$ws.Range("E48").Value = '  -4.64%  '
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E50").Value = '  -5.16%  '
$ws.Range("D51").Value = '3.204.37'
$ws.Range("E51").Value = '  -2.95%  '
